# Add a new "Italy" sheet to the workbook, mirroring the Excel UI flow used
# by the author: select the "Swiss" sheet, select all its cells (Ctrl+A),
# then Move-or-Copy -> "Create a copy" placed after the last existing sheet
# ("Portugal"). Afterwards rename the copy to "Italy" and localize its two
# market-specific cells.

$wb = $excel.ActiveWorkbook

$swiss = $wb.Worksheets.Item("Swiss")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Leave Swiss with a "select all" view state, as seen in the final workbook.
$swiss.Select()
$swiss.Range("A1:XFD1048576").Select() | Out-Null

# Create a copy of Swiss, inserted immediately after the last sheet.
$swiss.Copy($null, $lastSheet)

# The copy becomes the new last sheet and the active sheet/tab.
$italy = $wb.Worksheets.Item($wb.Worksheets.Count)
$italy.Name = "Italy"

# Localize the market name and part number for Italy.
$italy.Range("B2").Value = "Italy Market"
$italy.Range("B4").Value = "NGC-3145/T2155 "

# Final selection left on B4.
$italy.Range("B4").Select() | Out-Null
